$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1423601482268282
$ws.Range("D2").Value = 0.1726500903206158
$ws.Range("E2").Value = 0.1544409445399424
$ws.Range("F2").Value = 1.493400524436439
$ws.Range("G2").Value = 0.002465178597842065
$ws.Range("J2").Value = 0.1777225435731253
$ws.Range("K2").Value = 0.5777082146119312
$ws.Range("M2").Value = 0.2712404879345058
$ws.Range("N2").Value = 1.767832894244421
$ws.Range("O2").Value = 3.647913990833814
$ws.Range("B3").Value = 0.1329715801179958
$ws.Range("D3").Value = 0.1704119950685765
$ws.Range("E3").Value = 0.1538308296408708
$ws.Range("F3").Value = 1.494260510568978
$ws.Range("G3").Value = 0.00246784355435104
$ws.Range("J3").Value = 0.1779747742970024
$ws.Range("K3").Value = 0.5143002711891995
$ws.Range("M3").Value = 0.2556193327603751
$ws.Range("N3").Value = 1.785784428892363
$ws.Range("O3").Value = 3.656710742456653
$ws.Range("B4").Value = 0.127276348784946
$ws.Range("D4").Value = 0.169103000173223
$ws.Range("E4").Value = 0.1535225642610492
$ws.Range("F4").Value = 1.495550975367934
$ws.Range("G4").Value = 0.002469568523068895
$ws.Range("J4").Value = 0.1782117021526659
$ws.Range("K4").Value = 0.4753771775224038
$ws.Range("M4").Value = 0.2461189481287249
$ws.Range("N4").Value = 1.797403214017047
$ws.Range("O4").Value = 3.664046089926416
$ws.Range("B5").Value = 0.1249731139488688
$ws.Range("D5").Value = 0.1685860425347983
$ws.Range("E5").Value = 0.1534136605488854
$ws.Range("F5").Value = 1.496268683203184
$ws.Range("G5").Value = 0.002470293826081594
$ws.Range("J5").Value = 0.1783289023560357
$ws.Range("K5").Value = 0.4595189801710546
$ws.Range("M5").Value = 0.2422705947563415
$ws.Range("N5").Value = 1.802288006364915
$ws.Range("O5").Value = 3.667521733463019
$ws.Range("B6").Value = 0.1245917326544657
$ws.Range("D6").Value = 0.1685011988877534
$ws.Range("E6").Value = 0.1533965878117485
$ws.Range("F6").Value = 1.49639944764855
$ws.Range("G6").Value = 0.002470415614966987
$ws.Range("J6").Value = 0.178349611060078
$ws.Range("K6").Value = 0.4568859618524925
$ws.Range("M6").Value = 0.2416329824342043
$ws.Range("N6").Value = 1.803108186331915
$ws.Range("O6").Value = 3.668128243106082
$ws.Range("B7").Value = 0.1272452149932946
$ws.Range("D7").Value = 0.1690959615262315
$ws.Range("E7").Value = 0.1535210278168684
$ws.Range("F7").Value = 1.49555987775279
$ws.Range("G7").Value = 0.002469578214127756
$ws.Range("J7").Value = 0.1782131991220837
$ws.Range("K7").Value = 0.4751632939151307
$ws.Range("M7").Value = 0.2460669539287537
$ws.Range("N7").Value = 1.797468484441767
$ws.Range("O7").Value = 3.664090994028612
$ws.Range("B8").Value = 0.1391086675352255
$ws.Range("D8").Value = 0.1718649041167595
$ws.Range("E8").Value = 0.1542168259027576
$ws.Range("F8").Value = 1.493538869803032
$ws.Range("G8").Value = 0.002466079110395297
$ws.Range("J8").Value = 0.177792491138689
$ws.Range("K8").Value = 0.555843648427583
$ws.Range("M8").Value = 0.2658355456946211
$ws.Range("N8").Value = 1.773898820792159
$ws.Range("O8").Value = 3.650545729134222
$ws.Range("B9").Value = 0.1629172146540441
$ws.Range("D9").Value = 0.1778094095472937
$ws.Range("E9").Value = 0.1561064895344408
$ws.Range("F9").Value = 1.49562197067975
$ws.Range("G9").Value = 0.002459917948239054
$ws.Range("J9").Value = 0.1776180250091457
$ws.Range("K9").Value = 0.7141042124139858
$ws.Range("M9").Value = 0.3053162444351756
$ws.Range("N9").Value = 1.732408432586514
$ws.Range("O9").Value = 3.639329828638154
$ws.Range("B10").Value = 0.1807347358249984
$ws.Range("D10").Value = 0.1824874803157854
$ws.Range("E10").Value = 0.1578136859921599
$ws.Range("F10").Value = 1.500836636091151
$ws.Range("G10").Value = 0.00245581419378977
$ws.Range("J10").Value = 0.1778859461848228
$ws.Range("K10").Value = 0.8303779417543637
$ws.Range("M10").Value = 0.334750635999761
$ws.Range("N10").Value = 1.704803713151293
$ws.Range("O10").Value = 3.640450615427909
$ws.Range("B11").Value = 0.1889097726549522
$ws.Range("D11").Value = 0.1846824933857505
$ws.Range("E11").Value = 0.1586593135970418
$ws.Range("F11").Value = 1.504008724631632
$ws.Range("G11").Value = 0.002454038207542248
$ws.Range("J11").Value = 0.1780937551819974
$ws.Range("K11").Value = 0.8832683963954082
$ws.Range("M11").Value = 0.3482325790223371
$ws.Range("N11").Value = 1.692869544115929
$ws.Range("O11").Value = 3.642994770935417
$ws.Range("B12").Value = 0.1920153271915979
$ws.Range("D12").Value = 0.1855232482294582
$ws.Range("E12").Value = 0.1589894267470449
$ws.Range("F12").Value = 1.505324875135017
$ws.Range("G12").Value = 0.002453378680975258
$ws.Range("J12").Value = 0.1781847922989215
$ws.Range("K12").Value = 0.9032954795119394
$ws.Range("M12").Value = 0.3533508904681284
$ws.Range("N12").Value = 1.688439976948654
$ws.Range("O12").Value = 3.644250747531487
$ws.Range("B13").Value = 0.1913460555344244
$ws.Range("D13").Value = 0.1853417530269752
$ws.Range("E13").Value = 0.1589178915197067
$ws.Range("F13").Value = 1.505036307394079
$ws.Range("G13").Value = 0.002453520144677902
$ws.Range("J13").Value = 0.1781646370045564
$ws.Range("K13").Value = 0.8989823654493989
$ws.Range("M13").Value = 0.3522479966738885
$ws.Range("N13").Value = 1.68938997644365
$ws.Range("O13").Value = 3.643967238336046
$ws.Range("B14").Value = 0.1891650722420763
$ws.Range("D14").Value = 0.1847514718006664
$ws.Range("E14").Value = 0.158686274131
$ws.Range("F14").Value = 1.504114701808973
$ws.Range("G14").Value = 0.002453983687559226
$ws.Range("J14").Value = 0.1781009975313594
$ws.Range("K14").Value = 0.8849160697340608
$ws.Range("M14").Value = 0.3486534066360463
$ws.Range("N14").Value = 1.692503324323305
$ws.Range("O14").Value = 3.643092237375896
$ws.Range("B15").Value = 0.1878304336619721
$ws.Range("D15").Value = 0.1843911491311729
$ws.Range("E15").Value = 0.1585456891156767
$ws.Range("F15").Value = 1.503565159300948
$ws.Range("G15").Value = 0.002454269312236719
$ws.Range("J15").Value = 0.1780636237512851
$ws.Range("K15").Value = 0.8762998491260987
$ws.Range("M15").Value = 0.3464533009863331
$ws.Range("N15").Value = 1.694422014267428
$ws.Range("O15").Value = 3.642594374865951
$ws.Range("B16").Value = 0.1802018657913322
$ws.Range("D16").Value = 0.1823453716159804
$ws.Range("E16").Value = 0.1577598082706153
$ws.Range("F16").Value = 1.500645422338067
$ws.Range("G16").Value = 0.002455932081182736
$ws.Range("J16").Value = 0.177874093463366
$ws.Range("K16").Value = 0.8269212839546469
$ws.Range("M16").Value = 0.33387139102134
$ws.Range("N16").Value = 1.70559618310887
$ws.Range("O16").Value = 3.640325286817216
$ws.Range("B17").Value = 0.1755397211630623
$ws.Range("D17").Value = 0.181107444313426
$ws.Range("E17").Value = 0.1572953495615934
$ws.Range("F17").Value = 1.499059080098888
$ws.Range("G17").Value = 0.002456975357161711
$ws.Range("J17").Value = 0.1777798228243412
$ws.Range("K17").Value = 0.796627679754522
$ws.Range("M17").Value = 0.3261762053441686
$ws.Range("N17").Value = 1.712610833384208
$ws.Range("O17").Value = 3.639454329984687
$ws.Range("B18").Value = 0.1728647558710463
$ws.Range("D18").Value = 0.1804017259646855
$ws.Range("E18").Value = 0.1570347031658734
$ws.Range("F18").Value = 1.498221954015449
$ws.Range("G18").Value = 0.002457583975140147
$ws.Range("J18").Value = 0.1777336897505037
$ws.Range("K18").Value = 0.7792033694684619
$ws.Range("M18").Value = 0.3217588210503877
$ws.Range("N18").Value = 1.716704147388697
$ws.Range("O18").Value = 3.63914485517401
$ws.Range("B19").Value = 0.1719601938322342
$ws.Range("D19").Value = 0.1801638669543877
$ws.Range("E19").Value = 0.1569475697817069
$ws.Range("F19").Value = 1.497951452011662
$ws.Range("G19").Value = 0.002457791513602297
$ws.Range("J19").Value = 0.1777194595369096
$ws.Range("K19").Value = 0.7733037866576922
$ws.Range("M19").Value = 0.3202646695649989
$ws.Range("N19").Value = 1.71810015014815
$ws.Range("O19").Value = 3.639072957702808
$ws.Range("B20").Value = 0.1760353350656914
$ws.Range("D20").Value = 0.1812385718129264
$ws.Range("E20").Value = 0.1573441197005714
$ws.Range("F20").Value = 1.499220156637605
$ws.Range("G20").Value = 0.002456863413743624
$ws.Range("J20").Value = 0.1777890210125221
$ws.Range("K20").Value = 0.7998525163344254
$ws.Range("M20").Value = 0.326994474153679
$ws.Range("N20").Value = 1.711858039449719
$ws.Range("O20").Value = 3.639527227503834
$ws.Range("B21").Value = 0.1898054141792471
$ws.Range("D21").Value = 0.1849245930846166
$ws.Range("E21").Value = 0.1587540375528178
$ws.Range("F21").Value = 1.504382280628434
$ws.Range("G21").Value = 0.002453847181404625
$ws.Range("J21").Value = 0.1781193550631315
$ws.Range("K21").Value = 0.889047726859701
$ws.Range("M21").Value = 0.3497088742545031
$ws.Range("N21").Value = 1.691586425618418
$ws.Range("O21").Value = 3.643341306133948
$ws.Range("B22").Value = 0.1988622576007515
$ws.Range("D22").Value = 0.187389257366334
$ws.Range("E22").Value = 0.1597331455632833
$ws.Range("F22").Value = 1.508426004455842
$ws.Range("G22").Value = 0.00245195165211252
$ws.Range("J22").Value = 0.1784071941204957
$ws.Range("K22").Value = 0.9473334464460379
$ws.Range("M22").Value = 0.3646296601727883
$ws.Range("N22").Value = 1.678860270711187
$ws.Range("O22").Value = 3.647539305031927
$ws.Range("B23").Value = 0.1940232682266156
$ws.Range("D23").Value = 0.1860687534355776
$ws.Range("E23").Value = 0.1592053132326932
$ws.Range("F23").Value = 1.506206515403804
$ws.Range("G23").Value = 0.002452956420085647
$ws.Range("J23").Value = 0.1782469899437586
$ws.Range("K23").Value = 0.9162263517341103
$ws.Range("M23").Value = 0.3566593201075747
$ws.Range("N23").Value = 1.685604646089143
$ws.Range("O23").Value = 3.645142711169228
$ws.Range("B24").Value = 0.1758112512995069
$ws.Range("D24").Value = 0.181179270430718
$ws.Range("E24").Value = 0.1573220508546207
$ws.Range("F24").Value = 1.49914710067101
$ws.Range("G24").Value = 0.002456913996065864
$ws.Range("J24").Value = 0.1777848373914068
$ws.Range("K24").Value = 0.7983945929190668
$ws.Range("M24").Value = 0.3266245139974586
$ws.Range("N24").Value = 1.712198189137837
$ws.Range("O24").Value = 3.639493674799837
$ws.Range("B25").Value = 0.1564187622222448
$ws.Range("D25").Value = 0.1761464886592989
$ws.Range("E25").Value = 0.1555391948681049
$ws.Range("F25").Value = 1.494411322972084
$ws.Range("G25").Value = 0.002461510144612586
$ws.Range("J25").Value = 0.1775956281264968
$ws.Range("K25").Value = 0.6712885173101597
$ws.Range("M25").Value = 0.2945599664045488
$ws.Range("N25").Value = 1.743126875084027
$ws.Range("O25").Value = 3.64072052343235
